# =====================================================================
# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet named "2022-Q1" right before the "总计"
#    (grand-total) sheet, and populate it with the per-fund holding
#    detail for the 2022-Q1 quarter (mirrors the layout of the other
#    quarterly sheets: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#    持有市值(亿元)/仓位排名).
# 2. Insert a new row at the top of the "总计" (grand-total) sheet's
#    data for the 2022-Q1 quarter, pushing the older quarters down by
#    one row.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: make the header row (B..H or B..D) look like the other sheets
# -- bold font, thin border, centered/top aligned -- by copying the
# formatting straight from an existing header cell on another sheet
# (so the style gets de-duplicated onto the very same style index
# instead of minting a new one).
# ---------------------------------------------------------------------
function Copy-HeaderStyle($srcCell, $dstRange) {
    $srcCell.Copy()
    $dstRange.PasteSpecial(-4122)  # xlPasteFormats
}

$refSheet = $wb.Worksheets.Item("2021-Q4")
$refHeaderCell = $refSheet.Range("B1")
$refIndexCell = $refSheet.Range("A2")

# ---------------------------------------------------------------------
# Step 1: new "2022-Q1" worksheet, inserted immediately before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

$q1Headers = @('基金代码','基金名称','基金规模','股票总仓位','仓位占比','持有市值(亿元)','仓位排名')
for ($col = 2; $col -le 8; $col++) {
    $cell = $q1.Cells.Item(1, $col)
    $cell.Value = $q1Headers[$col - 2]
}
Copy-HeaderStyle $refHeaderCell $q1.Range("B1:H1")

$q1Data = @(
        @('003634','嘉实农业产业股票','28.18','94.29','7.63','2.1501',6),
        @('003751','万家瑞隆混合','27.84','86.40','6.79','1.8903',3),
        @('519196','万家新兴蓝筹灵活配置混合','21.26','80.70','6.76','1.4372',3),
        @('161810','银华内需精选混合(LOF)','25.59','94.71','5.13','1.3128',8),
        @('009199','万家价值优势一年持有期混合','13.70','89.53','8.67','1.1878',4),
        @('519195','万家品质生活灵活配置混合','17.66','79.87','6.53','1.1532',4),
        @('161912','万家社会责任18个月定期开放混合（LOF）A','13.56','88.11','8.20','1.1119',5),
        @('012412','汇泉策略优选混合型证券投资基金','23.52','70.35','4.67','1.0984',1),
        @('519181','万家和谐增长混合','11.51','88.62','9.31','1.0716',3),
        @('005094','万家臻选混合','13.43','73.12','5.06','0.6796',8),
        @('005106','银华农业产业股票','13.24','93.41','5.07','0.6713',6),
        @('001579','国泰大农业股票','11.72','90.32','3.31','0.3879',8),
        @('014107','博时品质生活混合A','4.11','69.12','6.23','0.2561',1),
        @('162208','泰达宏利首选企业股票','5.67','92.46','3.35','0.1899',8),
        @('001277','博时国企改革主题股票','2.77','89.22','6.64','0.1839',1),
        @('519767','交银施罗德科技创新灵活配置混合','4.18','90.38','4.01','0.1676',4),
        @('008234','光大保德信消费主题股票','2.26','91.66','5.58','0.1261',5),
        @('002535','中银鑫利灵活配置混合A','6.75','20.93','1.38','0.0932',1),
        @('006952','中银景元回报混合','3.58','33.73','1.68','0.0601',4),
        @('163823','中银稳健策略灵活配置混合','2.55','48.41','2.23','0.0569',5),
        @('008773','中银景泰回报混合','4.83','25.44','1.06','0.0512',7),
        @('002252','融通成长30灵活配置混合','1.70','78.26','2.61','0.0444',10),
        @('161913','万家社会责任18个月定期开放混合（LOF）C','0.44','88.11','8.20','0.0361',5),
        @('002536','中银鑫利灵活配置混合C','2.30','20.93','1.38','0.0317',1),
        @('002288','中银稳进策略灵活配置混合','0.70','66.42','4.29','0.0300',2),
        @('014108','博时品质生活混合C','0.14','69.12','6.23','0.0087',1),
        @('008112','中泰中证500指数增强A','0.61','92.46','0.89','0.0054',7),
        @('008113','中泰中证500指数增强C','0.46','92.46','0.89','0.0041',7),
        @('001614','东方区域发展混合','0.06','94.77','3.93','0.0024',10)
    )

for ($i = 0; $i -lt $q1Data.Count; $i++) {
    $r = $i + 2
    $row = $q1Data[$i]
    $q1.Cells.Item($r, 1).Value = $i
    $q1.Cells.Item($r, 2).Value = "'" + $row[0]
    $q1.Cells.Item($r, 3).Value = "'" + $row[1]
    $q1.Cells.Item($r, 4).Value = "'" + $row[2]
    $q1.Cells.Item($r, 5).Value = "'" + $row[3]
    $q1.Cells.Item($r, 6).Value = "'" + $row[4]
    $q1.Cells.Item($r, 7).Value = "'" + $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
}
Copy-HeaderStyle $refIndexCell $q1.Range("A2:A" + (1 + $q1Data.Count))

$q1.Range("A1").Select()

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q1 summary row at the top of "总计"'s data
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Range("B2").Value = "'2022-Q1"
$total.Range("C2").Value = 29
$total.Range("D2").Value = 15.5

$refTotalIndexCell = $total.Range("A3")
$refTotalIndexCell.Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("A2").Value = 0

$total.Range("A1").Select()
